# Add the "Chosen Interventions" sheet and populate it with a curated
# subset of rows copied from the "Interventions" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create the new sheet after the existing "Interventions" sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Chosen Interventions"

# xlPasteValues / xlPasteFormats paste-special constants
$xlPasteValues = -4163
$xlPasteFormats = -4122

# Header row (row 17 on Interventions -> row 1 on Chosen Interventions)
$ws1.Range("A17:D17").Copy() | Out-Null
$ws2.Range("A1:D1").PasteSpecial($xlPasteValues) | Out-Null
$ws1.Range("A17:D17").Copy() | Out-Null
$ws2.Range("A1:D1").PasteSpecial($xlPasteFormats) | Out-Null

# Chosen interventions, in the order they appear on the new sheet, and the
# source row on the "Interventions" sheet they were copied from:
#   row 2 <- 31 (Safety Campaigns)
#   row 3 <- 21 (Annual Health Check-ups)
#   row 4 <- 37 (Discounted Gym Memberships)
#   row 5 <- 24 (Weight Management Programs)
#   row 6 <- 54 (Cancer Prevention Initiatives)
#   row 7 <- 33 (Heart Health Screenings)
$rowMap = @{ 2 = 31; 3 = 21; 4 = 37; 5 = 24; 6 = 54; 7 = 33 }

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]

    $ws1.Range("A$($srcRow):D$($srcRow)").Copy() | Out-Null
    $ws2.Range("A$($destRow):D$($destRow)").PasteSpecial($xlPasteValues) | Out-Null

    $ws1.Range("A$($srcRow):D$($srcRow)").Copy() | Out-Null
    $ws2.Range("A$($destRow):D$($destRow)").PasteSpecial($xlPasteFormats) | Out-Null
}

# --- Column widths on the new sheet (approximate bestFit values) ---
$ws2.Columns.Item(1).ColumnWidth = 28.285714285714285
$ws2.Columns.Item(2).ColumnWidth = 69
$ws2.Columns.Item(3).ColumnWidth = 32.285714285714285
$ws2.Columns.Item(4).ColumnWidth = 23.428571428571427

# --- Selection / view state ---
$ws1.Range("A33:D33").Select() | Out-Null
$ws2.Range("C13").Select() | Out-Null

# Make the new sheet the active tab (also clears tabSelected on sheet 1)
$ws2.Activate() | Out-Null
